# NIT-9009306652.xlsx update:
#  - Refresh "Cant. Trabajadores" (E11) and "Cant. Periodos" (F13) totals
#  - Collapse the three worker/period detail rows down to a single row
#    (keep only period 1701 / 27578), removing the 1909 and 1908 rows
#  - Nudge the logo image a bit to the left

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Totals at the top of the statement
$ws.Range("E11").Value = 27578
$ws.Range("F13").Value = 1

# Keep a single detail row (period 1701 / value 27578) and drop the
# other two period rows (1909, 1908) that are no longer needed
$ws.Range("E16").Value = "1701"
$ws.Range("F16").Value = 27578
$ws.Rows("17:18").Delete()

# Move the company logo picture slightly to the left.
# (Shape.Left reads back rounded to 2 decimals, which is not precise
# enough to land on an exact EMU offset inside column B, so rebuild the
# current Left from the column-A width - which Range/Column Width
# reports at full precision - plus the image's original fractional
# offset into column B.)
$emuPerPoint = 12700.0
$colAWidthPt = $ws.Columns.Item(1).Width
$currentColOffEmu = 667900
$currentLeftPt = $colAWidthPt + ($currentColOffEmu / $emuPerPoint)
$shp = $ws.Shapes.Item(1)
$shp.Left = $currentLeftPt - (241300 / $emuPerPoint)
